$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '22.401.76'
$ws.Range("E2").Value = '  -4.65%  '
$ws.Range("D3").Value = '1.571.52'
$ws.Range("E3").Value = '  -4.65%  '
$ws.Range("E4").Value = '  -0.02%  '
Set-TextValue "D5" '1.001'
$ws.Range("E5").Value = '  -0.03%  '
Set-TextValue "D6" '290.94'
$ws.Range("E6").Value = '  -2.91%  '
Set-TextValue "D7" '0.3681'
$ws.Range("E7").Value = '  -2.88%  '
Set-TextValue "D8" '49.55'
$ws.Range("E8").Value = '  -2.12%  '
Set-TextValue "D9" '0.3371'
$ws.Range("E9").Value = '  -5.30%  '
$ws.Range("E10").Value = '  -4.64%  '
Set-TextValue "D11" '0.07561'
$ws.Range("E11").Value = '  -6.64%  '
Set-TextValue "D12" '1.001'
$ws.Range("E12").Value = '  -0.07%  '
Set-TextValue "D13" '21.09'
$ws.Range("E13").Value = '  -4.36%  '
Set-TextValue "D14" '6.056'
$ws.Range("E14").Value = '  -5.56%  '
$ws.Range("E15").Value = '  -7.44%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.574.35'
$ws.Range("E16").Value = '  -4.89%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D17" '0.00001137'
$ws.Range("E17").Value = '  -5.35%  '
Set-TextValue "D18" '89.28'
$ws.Range("E18").Value = '  -8.14%  '
Set-TextValue "D19" '0.06704'
$ws.Range("E19").Value = '  -3.96%  '
$ws.Range("E20").Value = '  +0.04%  '
Set-TextValue "D21" '6.222'
$ws.Range("E21").Value = '  -8.04%  '
$ws.Range("E22").Value = '  -6.52%  '
Set-TextValue "D23" '11.94'
$ws.Range("E23").Value = '  -4.63%  '
$ws.Range("D24").Value = '22.407.48'
$ws.Range("E24").Value = '  -4.74%  '
$ws.Range("E25").Value = '  -2.87%  '
Set-TextValue "D26" '2.956'
$ws.Range("E26").Value = '  +1.85%  '
$ws.Range("E27").Value = '  -5.54%  '
Set-TextValue "D28" '145.89'
$ws.Range("E28").Value = '  -4.51%  '
Set-TextValue "D29" '4.918'
$ws.Range("E29").Value = '  -5.76%  '
Set-TextValue "D30" '125.10'
$ws.Range("E30").Value = '  -5.97%  '
$ws.Range("D31").Value = '1.749.08'
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("E32").Value = '  -10.16%  '
Set-TextValue "D33" '1.970'
$ws.Range("E33").Value = '  -7.97%  '
Set-TextValue "D34" '0.9841'
$ws.Range("E34").Value = '  -4.03%  '
Set-TextValue "D35" '10.37'
$ws.Range("E35").Value = '  -12.74%  '
Set-TextValue "D36" '0.08443'
$ws.Range("E36").Value = '  -3.35%  '
Set-TextValue "D37" '0.02542'
$ws.Range("E37").Value = '  -6.88%  '
Set-TextValue "D38" '0.2295'
$ws.Range("E38").Value = '  -6.42%  '
Set-TextValue "D39" '0.06495'
$ws.Range("E39").Value = '  -5.09%  '
Set-TextValue "D40" '5.473'
$ws.Range("E40").Value = '  -8.33%  '
Set-TextValue "D41" '11.74'
$ws.Range("E41").Value = '  -12.86%  '
Set-TextValue "D42" '1.250'
$ws.Range("E42").Value = '  -5.30%  '
Set-TextValue "D43" '0.6372'
$ws.Range("E43").Value = '  -7.91%  '
Set-TextValue "D44" '14.47'
$ws.Range("E44").Value = '  -7.78%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E46").Value = '  -6.92%  '
Set-TextValue "D47" '3.770'
$ws.Range("E47").Value = '  -4.00%  '
Set-TextValue "D48" '2.108'
$ws.Range("E48").Value = '  -7.05%  '
Set-TextValue "D49" '121.23'
$ws.Range("E49").Value = '  -5.39%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue "D50" '1.195'
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D51" '0.07275'
$ws.Range("E51").Value = '  -6.67%  '
